$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range('A1').Value = '更新日期：2025.01.31 13:17:43'

# Shift daily-task rows 310-360: each row now shows the data that used to be one row below it
$rowData = @(
    @(310, '维荻', '9-3', '> 由非助战维荻累计造成60000点伤害> 3星通关主题曲9-3标准实战环境；必须编入非助战维荻并上场，且至少使用2次双刃毒藤'),
    @(311, '戴菲恩', 'WD-6', '> 由非助战戴菲恩累计造成30歼灭数> 3星通关插曲WD-6；必须携带且部署非助战戴菲恩，且至少使用2次抢攻'),
    @(312, '刺玫', 'IC-2', '> 战斗中非助战刺玫累计使用荆藤庇荫10次> 3星通关别传IC-2；必须编入非助战刺玫并上场，且不编入其他医疗干员'),
    @(313, '赫德雷', 'IW-7', '> 由非助战赫德雷累计歼灭5个精英或领袖敌人> 3星通关别传IW-7；必须编入非助战赫德雷并上场，并使用赫德雷至少击败2名沉沙'),
    @(314, '深律', 'LE-4', '> 完成5次战斗；必须编入非助战深律并上场，且每次战斗至少释放1次沉音宁神> 3星通关别传LE-4；必须编入非助战深律并上场，且所有干员不被击倒'),
    @(315, '止颂', '7-11', '> 由非助战止颂累计歼灭10个精英或领袖敌人> 3星通关主题曲7-11；必须编入非助战止颂并上场，且使用止颂歼灭至少2名雇佣军萨卡兹战士'),
    @(316, '止颂', 'TW-5', '> 由非助战止颂累计造成120000点伤害> 3星通关别传TW-5；必须编入非助战止颂并上场，且使用止颂击败至少6名敌人'),
    @(317, '薇薇安娜', 'MN-3', '> 由非助战薇薇安娜累计歼灭10个精英或领袖敌人> 3星通关别传MN-3；必须编入非助战薇薇安娜并上场，且使用薇薇安娜歼灭“锈铜”奥尔默·英格拉'),
    @(318, '塑心', 'GA-7', '> 由非助战塑心累计造成75000点凋亡损伤> 3星通关别传GA-7；必须编入非助战塑心并上场，且塑心造成至少15000点凋亡损伤'),
    @(319, '哈洛德', '9-13', '> 战斗中非助战哈洛德累计使用重症优先8次> 3星通关主题曲9-13标准实战环境；必须编入非助战哈洛德并上场，且所有干员不被击倒'),
    @(320, '烈夏', '4-3', '> 由非助战烈夏累计造成30歼灭数> 3星通关主题曲4-3；必须编入非助战烈夏并上场，且不编入其他近卫干员'),
    @(321, '锏', 'BI-6', '> 完成5次战斗；必须编入非助战锏并上场，且每次战斗至少释放1次归于宁静> 3星通关别传BI-6；必须编入非助战锏并上场，且使用锏至少歼灭10个敌人'),
    @(322, '莱伊', 'S9-1', '> 由非助战莱伊累计歼灭10个精英或领袖敌人> 3星通关主题曲S9-1标准实战环境；必须编入非助战莱伊并上场，且使用莱伊使用至少2次“得见光芒”'),
    @(323, '万顷', '9-13', '> 战斗中非助战万顷累计使用支援号令·γ型10次> 3星通关主题曲9-13标准实战环境；必须编入非助战万顷并上场，且至少使用2次应东风'),
    @(324, '小满', '9-11', '> 由非助战小满累计造成60000点伤害> 3星通关主题曲9-11标准实战环境；必须编入非助战小满并上场，且至少使用2次乡音沉沉'),
    @(325, '左乐', 'WB-5', '> 由非助战左乐累计造成40歼灭数> 3星通关别传WB-5；必须编入非助战左乐并上场，且使用左乐至少歼灭8个敌人'),
    @(326, '左乐', 'RI-7', '> 使用非助战左乐累计造成100000点伤害> 3星通关别传RI-7；必须携带且部署非助战左乐，且至少释放3次佑序有炎'),
    @(327, '黍', '11-11', '> 完成5次战斗；必须编入非助战黍并上场，且每次战斗至少释放1次离离枯荣> 3星通关主题曲11-11标准实战环境；必须编入非助战黍并上场，且所有干员不被击倒'),
    @(328, '红隼', '11-18', '> 完成5次战斗；必须编入非助战红隼并上场，且每次战斗至少释放2次醉刃乱舞> 3星通关主题曲11-18标准实战环境；必须编入非助战红隼，且第一位部署红隼、红隼全程不撤退或被击倒'),
    @(329, '导火索', '3-1', '> 由非助战导火索累计造成80000点伤害> 3星通关主题曲3-1；必须编入非助战导火索并上场，且使用导火索至少歼灭30个敌人'),
    @(330, '双月', '3-7', '> 战斗中非助战双月累计使用全知者的战术10次> 3星通关主题曲3-7；必须编入非助战双月并上场，且至少使用2次全知者的战术'),
    @(331, '医生', '2-5', '> 战斗中非助战医生累计使用激素手枪8次> 3星通关主题曲2-5；必须编入非助战医生并上场，且不编入医疗干员'),
    @(332, '艾拉', 'DM-EX-1', '> 战斗中非助战艾拉累计部署雷鸣地雷30个> 3星通关插曲DM-EX-1；必须编入非助战艾拉并上场，且使用艾拉歼灭至少2名萨卡兹穿刺手'),
    @(333, '露托', '3-1', '> 战斗中非助战露托累计使用强磁防卫8次> 3星通关主题曲3-1；必须编入非助战露托并上场，且不编入其他重装干员'),
    @(334, '奥达', 'S3-6', '> 使用非助战奥达累计使用4次锻锤之力> 3星通关主题曲S3-6；必须携带且部署非助战奥达，其他成员仅可编入4名干员'),
    @(335, '阿罗玛', 'GT-HX-3', '> 完成5次战斗；必须编入非助战阿罗玛并上场，且每次战斗至少释放1次小心地滑> 3星通关别传GT-HX-3；必须编入非助战阿罗玛并上场，且使用阿罗玛至少歼灭20个敌人'),
    @(336, '阿斯卡纶', '11-6', '> 由非助战阿斯卡纶累计造成180000点伤害> 3星通关主题曲11-6标准实战环境；必须编入非助战阿斯卡纶并上场，且使用阿斯卡纶至少歼灭30个敌人'),
    @(337, '历阵锐枪芬', '4-2', '> 完成5次战斗；必须编入非助战历阵锐枪芬，并确定第一位部署的干员是历阵锐枪芬> 3星通关主题曲4-2；必须编入非助战历阵锐枪芬并上场，且不编入其他先锋干员'),
    @(338, '魔王', '14-5', '> 携带非助战魔王完成5次战斗，且每次战斗至少发动一次“编织重构现世”> 3星通关主题曲14-5标准实战环境；必须编入非助战魔王并上场，其他成员仅可编入近战位干员'),
    @(339, '逻各斯', '11-6', '> 由非助战逻各斯累计造成100000点伤害> 3星通关主题曲11-6标准实战环境，必须编入非助战逻各斯并上场，其他成员仅可编入4名干员'),
    @(340, '维什戴尔', 'DM-5', '> 由非助战维什戴尔累计造成120000伤害> 3星通关插曲DM-5；必须编入非助战维什戴尔并上场，且使用维什戴尔至少歼灭20个敌人'),
    @(341, '阿米娅(医疗)', '3-8', '> 战斗中非助战医疗职业的阿米娅累计使用慈悲愿景5次> 3星通关主题曲3-8；必须编入非助战医疗职业的阿米娅并上场，且不编入其他医疗干员'),
    @(342, '深巡', 'SN-5', '> 战斗中非助战深巡累计使用行动能力剥夺8次> 3星通关插曲SN-5，必须编入非助战深巡并上场，其他成员仅可编入4名干员'),
    @(343, '海霓', 'SV-4', '> 战斗中非助战海霓累计使用阻滞性显色剂8次> 3星通关插曲SV-4；必须编入非助战海霓并上场，且不编入其他辅助干员'),
    @(344, '乌尔比安', 'SV-6', '> 战斗中非助战乌尔比安累计使用必须开辟的通路10次> 3星通关插曲SV-6；必须编入非助战乌尔比安并上场，并使用乌尔比安至少击败2名囊海爬行者'),
    @(345, '渡桥', '3-1', '> 战斗中非助战渡桥累计使用承压功率8次> 3星通关主题曲3-1；必须编入非助战渡桥并上场，且至少使用3次承压功率'),
    @(346, '衡沙', 'DV-2', '> 战斗中累计召唤非助战衡沙的召唤物20回> 3星通关别传DV-2；必须编入非助战衡沙并上场，其他成员仅可编入4名干员'),
    @(347, '佩佩', '3-8', '> 由非助战佩佩累计造成40歼灭数> 3星通关主题曲3-8；必须编入非助战佩佩并上场，且使用佩佩歼灭碎骨'),
    @(348, '森西', '1-12', '> 战斗中非助战森西累计使用团体魔物大餐6次> 3星通关主题曲1-12；必须编入非助战森西并上场，且所有干员不被击倒'),
    @(349, '齐尔查克', '4-3', '> 战斗中非助战齐尔查克累计使用随机应变6次> 3星通关主题曲4-3；必须编入非助战齐尔查克并上场，其他成员不可编入先锋干员'),
    @(350, '莱欧斯', '2-4', '> 完成5次战斗；必须编入非助战莱欧斯并上场，且每次战斗至少释放1次威吓战法> 3星通关主题曲2-4；必须编入非助战莱欧斯并上场，并使用莱欧斯至少击败1名高阶术师'),
    @(351, '玛露西尔', '5-10', '> 由非助战玛露西尔累计造成100000点伤害> 3星通关主题曲5-10；必须编入非助战玛露西尔并上场，且使用玛露西尔至少歼灭10名敌人'),
    @(352, '凯瑟琳', '11-7', '> 使用非助战凯瑟琳累计部署15个“支援装置”> 3星通关主题曲11-7标准实战环境；必须编入非助战凯瑟琳并上场，且凯瑟琳使用至少2次“战火淬炼”'),
    @(353, '波卜', '4-8', '> 由非助战波卜累计造成30次灼燃爆发> 3星通关主题曲4-8；必须编入非助战波卜并上场，且波卜使用至少2次“此路不通”'),
    @(354, '维娜·维多利亚', '9-5', '> 由非助战维娜·维多利亚累计造成120000点伤害> 3星通关主题曲9-5标准实战环境；必须编入非助战维娜·维多利亚并上场，其他成员仅可编入5名干员'),
    @(355, '裁度', '5-8', '> 完成5次战斗；必须编入非助战裁度并部署至少2次，且使用裁度击败至少4名敌人> 3星通关主题曲5-8；必须编入非助战裁度并上场，且至少束缚12次敌人'),
    @(356, '弑君者', '4-4', '> 完成5次战斗；每次战斗至少部署3次非助战弑君者> 3星通关主题曲4-4；必须编入非助战弑君者并上场，且不编入其他特种干员'),
    @(357, '忍冬', 'S2-3', '> 由非助战忍冬累计造成80000点伤害> 3星通关主题曲S2-3；必须编入非助战忍冬并上场，且使用忍冬击败至少24名敌人'),
    @(358, '荒芜拉普兰德', 'IS-8', '> 由非助战荒芜拉普兰德累计造成150000点伤害> 3星通关别传IS-8；必须编入非助战荒芜拉普兰德并上场，且荒芜拉普兰德使用至少2次逐猎狂飙或终幕·浩劫'),
    @(359, '瑰盐', '4-6', '> 战斗中非助战瑰盐累计使用绝妙的长效药呀8次> 3星通关主题曲4-6；必须编入非助战瑰盐并上场，且至少使用1次绝妙的长效药呀'),
    @(360, '行箸', '3-2', '> 使用非助战行箸累计使用8次食不厌精> 3星通关主题曲3-2，编入非助战行箸且所有干员不能被击倒'),
)

foreach ($item in $rowData) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
    $ws.Cells.Item($r, 3).Value = $item[3]
}